$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "CanvasEnMain" sheet by duplicating "Capas de Profundidad"
#    (sheet 1) so it inherits the exact same cell styles / column layout,
#    then drop it in as the last tab (after "Colisiones").
# ---------------------------------------------------------------------------
$wsSource = $wb.Worksheets.Item("Capas de Profundidad")
$wsLast   = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsSource.Copy($null, $wsLast)

$wsNew = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNew.Name = "CanvasEnMain"

# ---------------------------------------------------------------------------
# 2. Fill in the new content. The order below controls the order new shared
#    strings are interned in, matching how the sheet was originally typed
#    (column B down, then column C down, for each row pair).
# ---------------------------------------------------------------------------
$wsNew.Range("B3").Value = "Canvas"
$wsNew.Range("B4").Value = "CanvasLayout"
$wsNew.Range("C3").Value = "Uso"
$wsNew.Range("C4").Value = "Pausa"
$wsNew.Range("B5").Value = "CanvasLayout2"
$wsNew.Range("C5").Value = "Mapa"
$wsNew.Range("B6").Value = "CanvasLayout3"
$wsNew.Range("C6").Value = "Dialogos"

# ---------------------------------------------------------------------------
# 3. Column widths on the new sheet.
# ---------------------------------------------------------------------------
$wsNew.Columns.Item(2).ColumnWidth = 23.666666666666668
$wsNew.Columns.Item(3).ColumnWidth = 40.666666666666664

# Selection left on the new sheet after editing it.
$wsNew.Range("C14").Select()

# ---------------------------------------------------------------------------
# 4. Touch up the selections left on the other two sheets.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Capas de Profundidad")
$ws1.Activate()
$ws1.Columns.Item(2).ColumnWidth = 24.666666666666668
$ws1.Columns.Item(3).ColumnWidth = 24.833333333333336
$ws1.Range("B3:C15").Select()

$ws2 = $wb.Worksheets.Item("Colisiones")
$ws2.Activate()
$ws2.Columns.Item(2).ColumnWidth = 33.893229166666664
$ws2.Columns.Item(3).ColumnWidth = 33.619791666666664
$ws2.Columns.Item(4).ColumnWidth = 33.709635416666664
$ws2.Range("B8").Select()

# ---------------------------------------------------------------------------
# 5. Finish with the new sheet active/selected, as in the target workbook.
# ---------------------------------------------------------------------------
$wsNew.Activate()
$wsNew.Range("C14").Select()
